# comments + minor refactor
#
# The single data column (currently A1:A10, "Points on Bench" readings
# stored as text/shared-strings) is shifted one column to the right and
# one row down, to make room for a header row/label in column A.
# A new header cell, "Points on Bench", is written into B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data from column A into column B.
$ws.Columns("A").Insert()

# Shift everything down one row so row 1 is free for the new header.
$ws.Rows("1").Insert()

# B1 is now blank after the row insert; pick up the same cell style
# (font/format) already used by the data cells below it before writing
# the header text, so the new cell matches the rest of the column.
$ws.Range("B2").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "Points on Bench"
